$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.006.92"
$ws.Range("E2").Value = "  +5.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.486.43"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.36"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.34"
$ws.Range("E6").Value = "  +4.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.523"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.30"
$ws.Range("E10").Value = "  +2.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.31"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.15"
$ws.Range("E14").Value = "  +3.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.868.10"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.497.92"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.843.27"
$ws.Range("E18").Value = "  +5.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.64"
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.49"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0938"
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.64"
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.41"
$ws.Range("E23").Value = "  +5.91%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "251.35"
$ws.Range("E24").Value = "  +3.78%  "
$ws.Range("E25").Value = "  +2.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.28"
$ws.Range("E26").Value = "  +4.39%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.86"
$ws.Range("E29").Value = "  +3.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.09"
$ws.Range("E30").Value = "  +4.60%  "
$ws.Range("E31").Value = "  +7.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.63"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.72"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0771"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.92"
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("E39").Value = "  +3.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "122.90"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.02"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.980.09"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.98"
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("E48").Value = "  +3.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.02"
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.39"
$ws.Range("E50").Value = "  +17.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.79"
$ws.Range("E51").Value = "  +5.46%  "
